$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.2
    "C2" = 0.5666666666666667
    "J2" = 0.01944444444444444
    "P2" = 0.1472222222222222
    "S2" = 0.06666666666666667
    "C3" = 0.04147465437788019
    "J3" = 0.05069124423963134
    "P3" = 0.7926267281105991
    "S3" = 0.1152073732718894
    "J4" = 0.06666666666666667
    "P4" = 0.6666666666666666
    "S4" = 0.2666666666666667
    "B6" = 0.0975609756097561
    "D6" = 0.01463414634146342
    "F6" = 0.06341463414634146
    "J6" = 0.2195121951219512
    "O6" = 0.00975609756097561
    "Q6" = 0.2195121951219512
    "R6" = 0.05365853658536585
    "S6" = 0.3219512195121951
    "B7" = 0.1764705882352941
    "D7" = 0.0160427807486631
    "F7" = 0.0213903743315508
    "J7" = 0.1871657754010695
    "O7" = 0.0213903743315508
    "Q7" = 0.1657754010695187
    "R7" = 0.0481283422459893
    "S7" = 0.3636363636363636
    "B8" = 0.1083172147001934
    "D8" = 0.01740812379110251
    "F8" = 0.04642166344294004
    "J8" = 0.1237911025145068
    "O8" = 0.01740812379110251
    "Q8" = 0.2147001934235977
    "R8" = 0.07156673114119923
    "S8" = 0.4003868471953578
    "B9" = 0.09090909090909091
    "D9" = 0.0303030303030303
    "E9" = 0.006060606060606061
    "F9" = 0.04242424242424243
    "J9" = 0.06666666666666667
    "O9" = 0.01818181818181818
    "Q9" = 0.2
    "R9" = 0.103030303030303
    "S9" = 0.4424242424242424
    "B10" = 0.1141226818830243
    "D10" = 0.02995720399429386
    "E10" = 0.0007132667617689016
    "F10" = 0.07203994293865906
    "J10" = 0.1176890156918688
    "O10" = 0.01069900142653352
    "Q10" = 0.2417974322396576
    "R10" = 0.05563480741797432
    "S10" = 0.3573466476462197
    "G11" = 0.2014388489208633
    "J11" = 0.1043165467625899
    "K11" = 0.2266187050359712
    "L11" = 0.4496402877697842
    "S11" = 0.01798561151079137
    "F12" = 0.007936507936507936
    "G12" = 0.8015873015873016
    "J12" = 0.119047619047619
    "K12" = 0.01587301587301587
    "L12" = 0.03174603174603174
    "S12" = 0.02380952380952381
    "F13" = 0.02127659574468085
    "G13" = 0.6382978723404256
    "J13" = 0.2978723404255319
    "S13" = 0.0425531914893617
    "F15" = 0.01463414634146342
    "H15" = 0.1512195121951219
    "I15" = 0.05853658536585366
    "J15" = 0.375609756097561
    "K15" = 0.07317073170731707
    "M15" = 0.00975609756097561
    "O15" = 0.05853658536585366
    "S15" = 0.2585365853658537
    "F16" = 0.003875968992248062
    "H16" = 0.2170542635658915
    "I16" = 0.05038759689922481
    "J16" = 0.4302325581395349
    "K16" = 0.1162790697674419
    "M16" = 0.01550387596899225
    "N16" = 0.007751937984496124
    "O16" = 0.02713178294573643
    "S16" = 0.1317829457364341
    "F17" = 0.01088929219600726
    "H17" = 0.1814882032667877
    "I17" = 0.0852994555353902
    "J17" = 0.4573502722323049
    "K17" = 0.06715063520871144
    "M17" = 0.0235934664246824
    "O17" = 0.05807622504537205
    "S17" = 0.1161524500907441
    "F18" = 0.01324503311258278
    "H18" = 0.1920529801324503
    "I18" = 0.05960264900662252
    "J18" = 0.4768211920529801
    "K18" = 0.06622516556291391
    "O18" = 0.07947019867549669
    "S18" = 0.1125827814569536
    "F19" = 0.008547008547008548
    "H19" = 0.2362082362082362
    "I19" = 0.06526806526806526
    "J19" = 0.3954933954933955
    "K19" = 0.09090909090909091
    "M19" = 0.02408702408702409
    "O19" = 0.06604506604506605
    "S19" = 0.1134421134421134
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
